$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @(4333,0,268,728,51,1,2014,13472,4,407,16061)
$data[3] = @(821,0,446,452,4,1,89,8028,8,1,22)
$data[4] = @(254,0,73,75,2,0,267,4932,2,0,0)
$data[5] = @(803,0,261,266,3,2,114,7432,2,0,0)
$data[6] = @(255,0,124,130,5,1,388,5203,2,0,0)
$data[7] = @(251,0,9,9,0,0,0,6667,0,0,0)
$data[8] = @(772,0,166,181,15,0,829,2651,0,0,0)
$data[9] = @(774,0,26,26,0,0,0,9545,0,0,0)
$data[10] = @(559,1,49,51,2,0,1623,22943,3,0,0)
$data[11] = @(878,0,159,160,2,0,843,66784,2,0,0)
$data[12] = @(755,0,23,23,0,0,0,18947,1,0,0)
$data[13] = @(447,1,203,252,55,1,4794,5832,7,0,0)
$data[14] = @(804,0,161,163,0,2,0,8734,2,0,0)
$data[15] = @(253,0,61,64,3,0,469,4754,0,0,0)
$data[16] = @(172,1,341,453,25,2,2888,8243,5,84,25)
$data[17] = @(45,0,62,69,6,0,870,4032,1,1,145)
$data[18] = @(212,0,511,584,30,2,515,5949,9,41,702)
$data[19] = @(1511,0,171,181,7,3,807,10238,2,0,0)
$data[20] = @(1371,0,63,92,29,0,6468,1207,0,0,0)
$data[21] = @(410,0,39,37,0,0,0,24452,0,0,0)
$data[22] = @(2,1,0,1,0,0,0,0,0,0,0)
$data[23] = @(813,0,9,9,0,0,0,6667,1,0,0)
$data[24] = @(934,0,28,28,0,0,0,15000,1,0,0)
$data[25] = @(57,0,57,69,3,0,435,3860,1,9,1304)
$data[26] = @(116,0,0,0,0,0,0,0,0,0,0)
$data[27] = @(172,0,30,32,2,0,625,2000,2,0,0)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $arr = New-Object "object[,]" 1,11
    for ($i = 0; $i -lt 11; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $rng = $ws.Range("B" + $r + ":L" + $r)
    $rng.Value = $arr
}

Write-Host "Updated rows 2-27 (columns B:L) per diff"
